$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.586.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.929.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4824'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08226'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.013'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.922.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.118'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.333'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06875'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.599.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.687'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.178'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.161.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.421'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.102'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.016'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09605'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.610'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.560'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.389'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06388'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02289'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5969'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.929'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1848'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.477'
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07503'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5568'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.978'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '118.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.438'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.41%  '
